$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.628.26"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.782.85"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  -2.99%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("E11").Value = "  +2.53%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0837"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.18%  "
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").Value = "3.221.21"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "2.791.30"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.928"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").Value = "51.577.66"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("E19").Value = "  +4.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.44%  "
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.165"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +13.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("E30").Value = "  +7.84%  "
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "51.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.53%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0452"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0831"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.21%  "
$ws.Range("D46").Value = "2.133.60"
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("E48").Value = "  +5.01%  "
$ws.Range("E49").Value = "  -5.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.903"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.52%  "
